# The commit adds one new weekly price-sample row for "Cebollín baby"
# (Agrícola del Norte S.A. de Arica) right after the existing header/first
# rows, at sheet row 22. Every row that used to be at 22..109 shifts down
# by one (to 23..110), and the sheet's dimension grows from R109 to R110.
#
# We reproduce this by inserting a new row 22 (which pushes all the old
# data down automatically, exactly like the diff shows) and then filling
# that new row 22 with the new sample's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 22; rows 22-109 shift down to 23-110.
$ws.Rows("22:22").Insert()

# Populate the newly inserted row 22 with the new data point.
$ws.Range("A22").Value = 1
$ws.Range("B22").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C22").Value = "Arica y Parinacota"
$ws.Range("D22").Value = 44859
$ws.Range("E22").Value = 15
$ws.Range("F22").Value = 100112038
$ws.Range("G22").Value = "Cebollín baby"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 250
$ws.Range("K22").Value = 1400
$ws.Range("L22").Value = 1500
$ws.Range("M22").Value = 1450
$ws.Range("N22").Value = "$/paquete 1,5 a 2 kilos"
$ws.Range("O22").Value = "Región de Arica y Parinacota"
$ws.Range("P22").Value = 725
$ws.Range("Q22").Value = 2
$ws.Range("R22").Value = "Hortaliza"
